$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Reporte de Formatos" (sheet1): report moves from Q3 2021 to Q4 2021
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Reporte de Formatos")

# Reporting period + validation/update dates
$ws1.Range("B8").Value = 44470
$ws1.Range("C8").Value = 44561
$ws1.Range("H8").Value = 44571
$ws1.Range("I8").Value = 44571

# Row 3 grows taller to accommodate the now-wrapped labels in A3/D3
$ws1.Rows(3).RowHeight = 42.75

# A3 / D3 pick up word-wrap (they previously had none)
$ws1.Range("A3").WrapText = $true
$ws1.Range("D3").WrapText = $true

# G3 loses its word-wrap
$ws1.Range("G3").WrapText = $false

# H3 / I3 (previously borderless) get the same thin box border as the rest
$ws1.Range("H3:I3").Borders.LineStyle = 1

# G2 keeps its header look but loses its bottom border
$ws1.Range("G2:I2").Borders.Item(9).LineStyle = -4142

# Selection left on A9 after the edits
$ws1.Activate()
$ws1.Range("A9").Select()

# ---------------------------------------------------------------------------
# Sheet "Tabla_397514" (sheet2 - Responsables de recibir)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Tabla_397514")
$ws2.Activate()
$ws2.Range("C10").Select()

# ---------------------------------------------------------------------------
# Sheet "Tabla_397515" (sheet3 - Responsables de administrar)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Tabla_397515")
$ws3.Columns("B").ColumnWidth = 21
$ws3.Range("E4").Value = "Secretario  Administrativo"
$ws3.Activate()
$ws3.Range("E6").Select()

# ---------------------------------------------------------------------------
# Sheet "Tabla_397516" (sheet4 - Responsables de ejercer)
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Tabla_397516")
$ws4.Columns("B").ColumnWidth = 29.7109375
$ws4.Activate()
$ws4.Range("A4:E4").Select()

$ws1.Activate()
